# buy_sale_deal.docx pattern fix
#
# The authoritative (textual) change in this commit is a template-tag typo
# fix: the "currency" placeholder was missing its leading "$", reading
# "{currency}" instead of the correct "${currency}" merge-field syntax used
# by every other placeholder in the document. Everything else in the
# upstream diff is Word's own background spell/grammar-checker re-marking
# already-unchanged text (<w:proofErr w:type="spellStart|spellEnd|gramStart|
# gramEnd"/>) - cosmetic proofing bookmarks with no effect on document text,
# and not something settable through the Word object model.
#
# So: locate the "{currency}" placeholder and insert the missing "$"
# immediately before the opening brace.

$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute("{currency}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if ($r.Find.Found) {
    # Collapse to the start of the match (just before the "{") and type the
    # missing "$" there, turning "{currency}" into the correct "${currency}".
    $r.Collapse(1)
    $r.InsertBefore("$")
    Write-Output "Fixed currency placeholder: inserted missing '$' before {currency}"
} else {
    Write-Output "WARNING: '{currency}' placeholder not found - no change applied"
}
